# Edits made during key documents generation
# Adds a new "Buck MOSFET Specifications - Low Side" block (rows 41-59) to the
# MOSFETs sheet, relabels / restyles the "R Junction Mounting Base" row of the
# Boost section into "R Junction Case" with a box border around the whole
# detail block, updates the R-junction derate values used in the Delta T
# calculations, and tidies up the view state (selection) on both sheets.

$wb = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("MOSFETs")
$ws2 = $wb.Worksheets.Item("Diodes")

# --- Boost MOSFET section -------------------------------------------------
# Put a thin box border around the whole "RDs ON ... R Junction Case" label
# column of the Boost block (A9:A18).
$ws.Range("A9:A18").Borders.LineStyle = 1

# Relabel "R Junction Mounting Base" -> "R Junction Case" and update its value.
$ws.Range("A18").Value = "R Junction Case"
$ws.Range("B18").Value = 1.1

# --- Buck MOSFET section (existing, rows 21-39) ---------------------------
# R junction derate value changed (label/style unchanged here).
$ws.Range("B38").Value = 1.1

# --- New "Buck MOSFET Specifications - Low Side" block (rows 41-59) -------
$ws.Range("A41:B41").Merge() | Out-Null
$ws.Range("A41").Value = "Buck MOSFET Specifications - Low Side"

$ws.Range("A42").Value = "Continuous Peak"
$ws.Range("B42").Value = 116

$ws.Range("A43").Value = "Max Inrush"
$ws.Range("B43").Value = 445

$ws.Range("A44").Value = "Irms"
$ws.Range("B44").Value = 70

$ws.Range("A45").Value = "Iavg"
$ws.Range("B45").Value = 11.13

$ws.Range("A46").Value = "Max Vds"
$ws.Range("B46").Value = 112.4

$ws.Range("A47").Value = "Cont Max Vds"
$ws.Range("B47").Value = 73

$ws.Range("A48").Value = "Potential Specs"

$ws.Range("A49").Value = "RDs ON"
$ws.Range("B49").Value = 0.0038999999999999998

$ws.Range("A50").Value = "trise"
$ws.Range("B50").Value = 0.000000059

$ws.Range("A51").Value = "tfall"
$ws.Range("B51").Value = 0.000000014

$ws.Range("A52").Value = "Power Calcs"

$ws.Range("A53").Value = "Pcond"
$ws.Range("B53").Formula = "=B44*B44*B49"

$ws.Range("A54").Value = "Pon"
$ws.Range("B54").Formula = "=0.5*B47*B44*(B50+B51)*B57"

$ws.Range("A55").Value = "Poff"
$ws.Range("B55").Formula = "=B54"

$ws.Range("A56").Value = "P Total"
$ws.Range("B56").Formula = "=SUM(B53:B55)"

$ws.Range("A57").Value = "fsw"
$ws.Range("B57").Value = 23000

$ws.Range("A58").Value = "R Junction Mounting Base"
$ws.Range("B58").Value = 0.7

$ws.Range("A59").Value = "Delta T"
$ws.Range("B59").Formula = "=B58*B56"

# Match the same number-formats used by the existing trise/tfall/Pon/Poff rows.
$ws.Range("B50").Style = $ws.Range("B30").Style
$ws.Range("B51").Style = $ws.Range("B31").Style
$ws.Range("B54").Style = $ws.Range("B34").Style
$ws.Range("B55").Style = $ws.Range("B35").Style

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.1015625
$ws.Columns.Item(5).ColumnWidth = 16.3125

# --- View state --------------------------------------------------------------
# Diodes sheet selection moves to A2 (without leaving it the active sheet).
$ws2.Range("A2").Select() | Out-Null

# MOSFETs stays the active/tab-selected sheet, scrolled down to the new block,
# with B52 selected.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B52").Select() | Out-Null
